# "before modify merge kline jump"
# Remove two stale kline rows (20151021, 20151116), shift the remaining
# volume (column B) values down by 2 to correct the pre-merge jump, and
# append the two new trailing rows that come after the merge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the obsolete rows for 20151021 (row 2) and 20151116 (row 3).
#    Deleting both in one shot avoids re-indexing issues.
$ws.Range("A2:A3").EntireRow.Delete()

# 2) Every remaining row's B value needs to drop by 2 (kline jump fix).
for ($r = 1; $r -le 56; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    $cell.Value2 = $current - 2
}

# 3) Append the two new rows that appear after the merge.
$ws.Cells.Item(57, 1).NumberFormat = "@"
$ws.Cells.Item(57, 1).Value2 = "20200428"
$ws.Cells.Item(57, 2).Value2 = 523.0

$ws.Cells.Item(58, 1).NumberFormat = "@"
$ws.Cells.Item(58, 1).Value2 = "20200603"
$ws.Cells.Item(58, 2).Value2 = 598.0
